$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 865-866; everything from the old row 865
# onward shifts down by two (old row 865 -> new row 867, ..., old row 934 -> new row 936).
$ws.Range("A865:A866").EntireRow.Insert()

# New row 865 data
$ws.Range("A865").Value = 4
$ws.Range("B865").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C865").Value = "Los Lagos"
$ws.Range("D865").Value = 45106
$ws.Range("E865").Value = 10
$ws.Range("F865").Value = 100112004
$ws.Range("G865").Value = "Cebolla"
$ws.Range("H865").Value = "Morada(o)"
$ws.Range("I865").Value = "1a (guarda)"
$ws.Range("J865").Value = 150
$ws.Range("K865").Value = 14000
$ws.Range("L865").Value = 14000
$ws.Range("M865").Value = 14000
$ws.Range("N865").Value = "`$/malla 18 kilos"
$ws.Range("O865").Value = "Región de O'Higgins"
$ws.Range("P865").Value = 778
$ws.Range("Q865").Value = 18
$ws.Range("R865").Value = "Hortaliza"

# New row 866 data
$ws.Range("A866").Value = 4
$ws.Range("B866").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C866").Value = "Los Lagos"
$ws.Range("D866").Value = 45106
$ws.Range("E866").Value = 10
$ws.Range("F866").Value = 100112004
$ws.Range("G866").Value = "Cebolla"
$ws.Range("H866").Value = "Sin especificar"
$ws.Range("I866").Value = "1a (guarda)"
$ws.Range("J866").Value = 500
$ws.Range("K866").Value = 12000
$ws.Range("L866").Value = 12000
$ws.Range("M866").Value = 12000
$ws.Range("N866").Value = "`$/malla 18 kilos"
$ws.Range("O866").Value = "Región de O'Higgins"
$ws.Range("P866").Value = 667
$ws.Range("Q866").Value = 18
$ws.Range("R866").Value = "Hortaliza"
